$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.991.12"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.849.34"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'309.32"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4776"
$ws.Range("E7").Value = "  +2.05%  "
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "'0.07231"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'0.9312"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "'19.79"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "'0.07737"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "1.858.84"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'5.350"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "'6.450"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "'89.08"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "'1.015"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'0.000008637"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "'1.012"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "27.006.92"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "'14.53"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "'1.933"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "'152.91"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'18.23"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").Value = "'2.015"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "'114.37"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "'4.963"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").Value = "'0.08860"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").Value = "'3.312"
$ws.Range("E31").Value = "  +4.76%  "
$ws.Range("D32").Value = "'1.183"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "'0.7416"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'4.506"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("D36").Value = "'1.110"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("D37").Value = "'0.01960"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").Value = "'0.05264"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("D41").Value = "'7.024"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "'0.1518"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").Value = "'8.255"
$ws.Range("E43").Value = "  +1.78%  "
$ws.Range("D44").Value = "'10.61"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("D45").Value = "'0.4754"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("D46").Value = "'1.012"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "'101.79"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "'1.608"
$ws.Range("D49").Value = "'65.84"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("D50").Value = "'0.06075"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "'0.8896"
$ws.Range("E51").Value = "  +4.12%  "
